{"js": "// Add a new work-log row at the end of the first (work log) table:\n// Date | Activity | Time\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst logTable = tables.items[0];\nlogTable.addRows(\n  \"End\",\n  1,\n  [\n    [\n      \"28/3/24\",\n      \"Worked on design report. Edited specifications and system overview. Started on ICT security risks and found two more references.\",\n      \"5\"\n    ]\n  ]\n);\nawait context.sync();\n", "ps1": "# Add a new work-log row at the end of the first (work log) table:\n# Date | Activity | Time\n$d = $word.ActiveDocument\n$logTable = $d.Tables(1)\n\n$newRow = $logTable.Rows.Add()\n$newRow.Cells(1).Range.Text = \"28/3/24\"\n$newRow.Cells(2).Range.Text = \"Worked on design report. Edited specifications and system overview. Started on ICT security risks and found two more references.\"\n$newRow.Cells(3).Range.Text = \"5\"\n"}
